# Weekly refresh of the Fruta/Hortaliza "Breva" sheet: the data columns
# (Fecha, Calidad, Volumen, Precio min/max/ponderado, Unidad, Origen,
# Precio $/Kg, Kg/unidad) are rotated across rows 2,3,5,6,7,10,11,12,13
# while the descriptive columns (A,B,C,E..K) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot of the "moving" columns for every affected row, keyed by row
# number, captured BEFORE any writes so the rotation reads consistent
# source data regardless of write order.
$rows = 2,3,5,6,7,10,11,12,13
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = [ordered]@{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Destination row -> source row the new values are pulled from.
$sourceFor = @{
    2  = 11
    3  = 5
    5  = 6
    6  = 7
    7  = 10
    10 = 12
    11 = 13
    12 = 2
    13 = 3
}

foreach ($dest in $rows) {
    $src = $snapshot[$sourceFor[$dest]]

    $ws.Cells.Item($dest, 4).Value  = $src.D
    $ws.Cells.Item($dest, 12).Value = $src.L
    $ws.Cells.Item($dest, 13).Value = $src.M
    $ws.Cells.Item($dest, 14).Value = $src.N
    $ws.Cells.Item($dest, 15).Value = $src.O
    $ws.Cells.Item($dest, 16).Value = $src.P
    $ws.Cells.Item($dest, 17).Value = $src.Q
    $ws.Cells.Item($dest, 18).Value = $src.R
    $ws.Cells.Item($dest, 19).Value = $src.S
    $ws.Cells.Item($dest, 20).Value = $src.T
}
